$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 57.14035266666667
$ws.Range("H2").Value = 171.421058
$ws.Range("I2").Value = 0.7274038390747541
$ws.Range("J2").Value = 0.7274038390747541
$ws.Range("M2").Value = 36.89194233333333
$ws.Range("N2").Value = 110.675827
$ws.Range("O2").Value = 0.3567095043190808
$ws.Range("P2").Value = 0.3567095043190809
$ws.Range("Q2").Value = 2108.018595484996
$ws.Range("R2").Value = 18972.16735936497
$ws.Range("S2").Value = 0.259471862876152
$ws.Range("T2").Value = 0.259471862876152

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 57.14035266666667
$ws.Range("H3").Value = 171.421058
$ws.Range("I3").Value = 0.7274038390747541
$ws.Range("J3").Value = 0.7274038390747541
$ws.Range("M3").Value = 42.68037399999999
$ws.Range("O3").Value = 0.4126780562577495
$ws.Range("P3").Value = 0.4126780562577496
$ws.Range("Q3").Value = 2438.771622305231
$ws.Range("R3").Value = 21948.94460074707
$ws.Range("S3").Value = 0.3001836024237944
$ws.Range("T3").Value = 0.3001836024237944

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 57.14035266666667
$ws.Range("H4").Value = 171.421058
$ws.Range("I4").Value = 0.7274038390747541
$ws.Range("J4").Value = 0.7274038390747541
$ws.Range("M4").Value = 23.85061433333334
$ws.Range("N4").Value = 71.551843
$ws.Range("O4").Value = 0.2306124394231696
$ws.Range("P4").Value = 0.2306124394231696
$ws.Range("Q4").Value = 1362.832514323322
$ws.Range("R4").Value = 12265.4926289099
$ws.Range("S4").Value = 0.1677483737748077
$ws.Range("T4").Value = 0.1677483737748078

# Row 5
$ws.Range("I5").Value = 0.08622113322131104
$ws.Range("J5").Value = 0.08622113322131104
$ws.Range("M5").Value = 36.89194233333333
$ws.Range("N5").Value = 110.675827
$ws.Range("O5").Value = 0.3567095043190808
$ws.Range("P5").Value = 0.3567095043190809
$ws.Range("Q5").Value = 249.8691131263525
$ws.Range("R5").Value = 2248.822018137173
$ws.Range("S5").Value = 0.03075589769320329
$ws.Range("T5").Value = 0.0307558976932033

# Row 6
$ws.Range("I6").Value = 0.08622113322131104
$ws.Range("J6").Value = 0.08622113322131104
$ws.Range("M6").Value = 42.68037399999999
$ws.Range("O6").Value = 0.4126780562577495
$ws.Range("P6").Value = 0.4126780562577496
$ws.Range("R6").Value = 2601.667429876877
$ws.Range("S6").Value = 0.03558156966611111
$ws.Range("T6").Value = 0.03558156966611112

# Row 7
$ws.Range("I7").Value = 0.08622113322131104
$ws.Range("J7").Value = 0.08622113322131104
$ws.Range("M7").Value = 23.85061433333334
$ws.Range("N7").Value = 71.551843
$ws.Range("O7").Value = 0.2306124394231696
$ws.Range("P7").Value = 0.2306124394231696
$ws.Range("S7").Value = 0.01988366586199663
$ws.Range("T7").Value = 0.01988366586199663

# Row 8
$ws.Range("I8").Value = 0.1863750277039348
$ws.Range("J8").Value = 0.1863750277039348
$ws.Range("M8").Value = 36.89194233333333
$ws.Range("N8").Value = 110.675827
$ws.Range("O8").Value = 0.3567095043190808
$ws.Range("P8").Value = 0.3567095043190809
$ws.Range("Q8").Value = 540.115411864839
$ws.Range("R8").Value = 4861.038706783551
$ws.Range("S8").Value = 0.06648174374972553
$ws.Range("T8").Value = 0.06648174374972556

# Row 9
$ws.Range("I9").Value = 0.1863750277039348
$ws.Range("J9").Value = 0.1863750277039348
$ws.Range("M9").Value = 42.68037399999999
$ws.Range("O9").Value = 0.4126780562577495
$ws.Range("P9").Value = 0.4126780562577496
$ws.Range("Q9").Value = 624.8607778161539
$ws.Range("R9").Value = 5623.747000345385
$ws.Range("S9").Value = 0.07691288416784402
$ws.Range("T9").Value = 0.07691288416784403

# Row 10
$ws.Range("I10").Value = 0.1863750277039348
$ws.Range("J10").Value = 0.1863750277039348
$ws.Range("M10").Value = 23.85061433333334
$ws.Range("N10").Value = 71.551843
$ws.Range("O10").Value = 0.2306124394231696
$ws.Range("P10").Value = 0.2306124394231696
$ws.Range("Q10").Value = 1362.832514323322
$ws.Range("R10").Value = 3142.658047314159
$ws.Range("S10").Value = 0.04298039978636521
$ws.Range("T10").Value = 0.04298039978636523
